$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record attendance for the most recent Journal session (column G) for every
# attendee (rows 3-8). Everyone attended this session, so each cell becomes 1.
# The "Total" (column B) and "Maximum" (row 9) formulas already in the sheet
# will recompute automatically from these new values.
$ws.Range("G3:G8").Value = 1

# Reflect where the user's selection ended up after entering the new data.
$ws.Range("F12").Select()
